$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2395.8235
$ws.Range("I40").Value = 3471.5
$ws.Range("J40").Value = 1809.091
$ws.Range("K40").Value = 3471.5
$ws.Range("L40").Value = 1809.091
$ws.Range("M40").Value = -3296.5
$ws.Range("N40").Value = -2159.091

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 42033.297
$ws.Range("I98").Value = 53847.57
$ws.Range("J98").Value = 683.3333
$ws.Range("K98").Value = 53847.57
$ws.Range("L98").Value = 683.3333
$ws.Range("M98").Value = -52349.57
$ws.Range("N98").Value = -3679.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3999.375
$ws.Range("I116").Value = 1499
$ws.Range("K116").Value = 1499
$ws.Range("M116").Value = 1943

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 42033.297
$ws.Range("I122").Value = 53847.57
$ws.Range("J122").Value = 683.3333
$ws.Range("K122").Value = 161542.71
$ws.Range("L122").Value = 2049.9999
$ws.Range("M122").Value = -159092.71
$ws.Range("N122").Value = -6949.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2305932.5
$ws.Range("I132").Value = 2917324
$ws.Range("J132").Value = 1456.8462
$ws.Range("K132").Value = 8751972
$ws.Range("L132").Value = 4370.5386
$ws.Range("M132").Value = -8749442
$ws.Range("N132").Value = -9430.5386

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1807.022
$ws.Range("I138").Value = 1445.55
$ws.Range("J138").Value = 2090.5293
$ws.Range("K138").Value = 4336.65
$ws.Range("L138").Value = 6271.5879
$ws.Range("M138").Value = 803.3500000000004
$ws.Range("N138").Value = -16551.5879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33334128
$ws.Range("I2").Value = 55556084
$ws.Range("K2").Value = 55556084
$ws.Range("M2").Value = -55555971

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1357.95
$ws.Range("I32").Value = 1306.0714
$ws.Range("J32").Value = 3900
$ws.Range("K32").Value = 1306.0714
$ws.Range("L32").Value = 3900
$ws.Range("M32").Value = -1019.0714
$ws.Range("N32").Value = -4474

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 13889804
$ws.Range("I45").Value = 30303618
$ws.Range("J45").Value = 1192.3077
$ws.Range("K45").Value = 30303618
$ws.Range("L45").Value = 1192.3077
$ws.Range("M45").Value = -30303241
$ws.Range("N45").Value = -1946.3077

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2085.238
$ws.Range("I61").Value = 1660.4
$ws.Range("J61").Value = 3147.3333
$ws.Range("K61").Value = 1660.4
$ws.Range("L61").Value = 3147.3333
$ws.Range("M61").Value = -1448.4
$ws.Range("N61").Value = -3571.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 950.7843
$ws.Range("I74").Value = 869.34784
$ws.Range("K74").Value = 869.34784
$ws.Range("M74").Value = 4.652159999999981

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 20000
$ws.Range("J75").Value = 20000
$ws.Range("L75").Value = 20000
$ws.Range("N75").Value = -21748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 950.7843
$ws.Range("I77").Value = 869.34784
$ws.Range("K77").Value = 4346.7392
$ws.Range("M77").Value = 21.26080000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H78").Value = 20000
$ws.Range("J78").Value = 20000
$ws.Range("L78").Value = 60000
$ws.Range("N78").Value = -68736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2309.9375
$ws.Range("I110").Value = 2374.7144
$ws.Range("J110").Value = 1856.5
$ws.Range("K110").Value = 2374.7144
$ws.Range("L110").Value = 1856.5
$ws.Range("M110").Value = -329.7143999999998
$ws.Range("N110").Value = -5946.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 33334128
$ws.Range("I116").Value = 55556084
$ws.Range("K116").Value = 55556084
$ws.Range("M116").Value = -55553790

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1511.9688
$ws.Range("I122").Value = 1498.1666
$ws.Range("J122").Value = 1553.375
$ws.Range("K122").Value = 4494.4998
$ws.Range("L122").Value = 4660.125
$ws.Range("M122").Value = -2044.4998
$ws.Range("N122").Value = -9560.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2085.238
$ws.Range("I136").Value = 1660.4
$ws.Range("J136").Value = 3147.3333
$ws.Range("K136").Value = 4981.200000000001
$ws.Range("L136").Value = 9441.999899999999
$ws.Range("M136").Value = -2431.200000000001
$ws.Range("N136").Value = -14541.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33334128
$ws.Range("I3").Value = 55556084
$ws.Range("K3").Value = 55556084
$ws.Range("M3").Value = -55555970

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 348.08694
$ws.Range("I22").Value = 348.08694
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 348.08694
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -175.08694
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 341.94446
$ws.Range("I80").Value = 428.5
$ws.Range("J80").Value = 317.2143
$ws.Range("K80").Value = 428.5
$ws.Range("L80").Value = 317.2143
$ws.Range("M80").Value = 569.5
$ws.Range("N80").Value = -2313.2143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 341.94446
$ws.Range("I83").Value = 428.5
$ws.Range("J83").Value = 317.2143
$ws.Range("K83").Value = 2142.5
$ws.Range("L83").Value = 1586.0715
$ws.Range("M83").Value = 2849.5
$ws.Range("N83").Value = -11570.0715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2615.1072
$ws.Range("I86").Value = 2295.9473
$ws.Range("J86").Value = 3288.889
$ws.Range("K86").Value = 2295.9473
$ws.Range("L86").Value = 3288.889
$ws.Range("M86").Value = -1172.9473
$ws.Range("N86").Value = -5534.889

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2615.1072
$ws.Range("I89").Value = 2295.9473
$ws.Range("J89").Value = 3288.889
$ws.Range("K89").Value = 11479.7365
$ws.Range("L89").Value = 16444.445
$ws.Range("M89").Value = -5863.736499999999
$ws.Range("N89").Value = -27676.445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1133.1666
$ws.Range("I99").Value = 1159.8
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1159.8
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 338.2
$ws.Range("N99").Value = -3996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1891.579
$ws.Range("I134").Value = 1587.9744
$ws.Range("J134").Value = 2549.389
$ws.Range("K134").Value = 4763.9232
$ws.Range("L134").Value = 7648.167
$ws.Range("M134").Value = -2228.9232
$ws.Range("N134").Value = -12718.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1035.4865
$ws.Range("I58").Value = 1053.3235
$ws.Range("K58").Value = 1053.3235
$ws.Range("M58").Value = -850.3235

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 861.9091
$ws.Range("I105").Value = 760
$ws.Range("K105").Value = 760
$ws.Range("M105").Value = 987

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 715.7727
$ws.Range("I107").Value = 670.6786
$ws.Range("J107").Value = 794.6875
$ws.Range("K107").Value = 670.6786
$ws.Range("L107").Value = 794.6875
$ws.Range("M107").Value = 1249.3214
$ws.Range("N107").Value = -4634.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 973.4091
$ws.Range("I122").Value = 885.41174
$ws.Range("J122").Value = 1272.6
$ws.Range("K122").Value = 2656.23522
$ws.Range("L122").Value = 3817.8
$ws.Range("M122").Value = -206.23522
$ws.Range("N122").Value = -8717.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1036.0714
$ws.Range("I134").Value = 1036.0714
$ws.Range("K134").Value = 3108.2142
$ws.Range("M134").Value = -573.2142000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1035.4865
$ws.Range("I136").Value = 1053.3235
$ws.Range("K136").Value = 3159.9705
$ws.Range("M136").Value = -609.9704999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 274.1875
$ws.Range("I2").Value = 357.75
$ws.Range("J2").Value = 23.5
$ws.Range("K2").Value = 2146.5
$ws.Range("L2").Value = 141
$ws.Range("M2").Value = -2033.5
$ws.Range("N2").Value = -367

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 37046044
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 6930.5
$ws.Range("I116").Value = 9229.166999999999
$ws.Range("J116").Value = 2333.1667
$ws.Range("K116").Value = 27687.501
$ws.Range("L116").Value = 6999.500100000001
$ws.Range("M116").Value = -24245.501
$ws.Range("N116").Value = -13883.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 38721956
$ws.Range("I70").Value = 50337308
$ws.Range("J70").Value = 4101
$ws.Range("K70").Value = 50337308
$ws.Range("L70").Value = 4101
$ws.Range("M70").Value = -50337038
$ws.Range("N70").Value = -4641

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 38721956
$ws.Range("I73").Value = 50337308
$ws.Range("J73").Value = 4101
$ws.Range("K73").Value = 50337308
$ws.Range("L73").Value = 4101
$ws.Range("M73").Value = -50336372
$ws.Range("N73").Value = -5973

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 383
$ws.Range("I107").Value = 352.4
$ws.Range("J107").Value = 444.2
$ws.Range("K107").Value = 352.4
$ws.Range("L107").Value = 444.2
$ws.Range("M107").Value = 1567.6
$ws.Range("N107").Value = -4284.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 15626132
$ws.Range("I113").Value = 41667400
$ws.Range("J113").Value = 1370
$ws.Range("K113").Value = 41667400
$ws.Range("L113").Value = 1370
$ws.Range("M113").Value = -41665230
$ws.Range("N113").Value = -5710

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1611.5454
$ws.Range("I126").Value = 1599.24
$ws.Range("K126").Value = 4797.72
$ws.Range("M126").Value = -2327.72

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1929.7778
$ws.Range("I40").Value = 1961.2693
$ws.Range("J40").Value = 1111
$ws.Range("K40").Value = 1961.2693
$ws.Range("L40").Value = 1111
$ws.Range("M40").Value = -1825.2693
$ws.Range("N40").Value = -1383

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 16669211
$ws.Range("I100").Value = 23811452
$ws.Range("J100").Value = 3981.3333
$ws.Range("K100").Value = 23811452
$ws.Range("L100").Value = 3981.3333
$ws.Range("M100").Value = -23810911
$ws.Range("N100").Value = -5063.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5338.769
$ws.Range("I122").Value = 7629.143
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 22887.429
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -20437.429
$ws.Range("N122").Value = -12900.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 16369411
$ws.Range("J107").Value = 33333694
$ws.Range("L107").Value = 100001082
$ws.Range("N107").Value = -100004922

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1253.9846
$ws.Range("I132").Value = 1116.1296
$ws.Range("J132").Value = 1930.7273
$ws.Range("K132").Value = 3348.3888
$ws.Range("L132").Value = 5792.1819
$ws.Range("M132").Value = -818.3887999999997
$ws.Range("N132").Value = -10852.1819
